$wb = $excel.ActiveWorkbook
$wsOut = $wb.Worksheets.Item("out")

# Shift existing data rows 2-7 down to rows 3-8 (bottom-up so we don't clobber
# a row before it has been read), then write the new first row of data.
for ($r = 7; $r -ge 2; $r--) {
    $a = $wsOut.Cells.Item($r, 1).Value()
    $b = $wsOut.Cells.Item($r, 2).Value()
    $c = $wsOut.Cells.Item($r, 3).Value()
    $d = $wsOut.Cells.Item($r, 4).Value()

    $wsOut.Cells.Item($r + 1, 1).Value = $a
    $wsOut.Cells.Item($r + 1, 2).Value = $b
    $wsOut.Cells.Item($r + 1, 3).Value = $c
    $wsOut.Cells.Item($r + 1, 4).Value = $d
}

$wsOut.Cells.Item(2, 1).Value = 20190812
$wsOut.Cells.Item(2, 2).Value = "keishi"
$wsOut.Cells.Item(2, 3).Value = 80
$wsOut.Cells.Item(2, 4).Value = "服务器首月"

# View / selection bookkeeping: "in" is no longer the active tab, "out" is.
$wsIn = $wb.Worksheets.Item("in")
$wsIn.Activate()
$excel.ActiveWindow.Zoom = 150

$wsOut.Activate()
$excel.ActiveWindow.Zoom = 150
$null = $wsOut.Range("E8").Select()
